$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 154
$ws.Cells.Item(2, 4).Value = 1.298701298701299
$ws.Cells.Item(4, 3).Value = 172
$ws.Cells.Item(4, 4).Value = 2.906976744186046
$ws.Cells.Item(5, 2).Value = 5
$ws.Cells.Item(5, 3).Value = 183
$ws.Cells.Item(5, 4).Value = 2.73224043715847
$ws.Cells.Item(6, 3).Value = 165
$ws.Cells.Item(6, 4).Value = 1.212121212121212
$ws.Cells.Item(7, 3).Value = 174
$ws.Cells.Item(7, 4).Value = 1.724137931034483
$ws.Cells.Item(8, 3).Value = 162
$ws.Cells.Item(8, 4).Value = 0.6172839506172839
$ws.Cells.Item(9, 2).Value = 6
$ws.Cells.Item(9, 3).Value = 134
$ws.Cells.Item(9, 4).Value = 4.477611940298507
$ws.Cells.Item(10, 3).Value = 116
$ws.Cells.Item(10, 4).Value = 4.310344827586207
$ws.Cells.Item(11, 3).Value = 109
$ws.Cells.Item(11, 4).Value = 2.752293577981652
$ws.Cells.Item(12, 3).Value = 98
$ws.Cells.Item(12, 4).Value = 2.040816326530612
$ws.Cells.Item(13, 3).Value = 116
$ws.Cells.Item(13, 4).Value = 2.586206896551724
$ws.Cells.Item(14, 2).Value = 6
$ws.Cells.Item(14, 3).Value = 126
$ws.Cells.Item(14, 4).Value = 4.761904761904762
$ws.Cells.Item(15, 3).Value = 230
$ws.Cells.Item(15, 4).Value = 4.782608695652174
$ws.Cells.Item(16, 3).Value = 224
$ws.Cells.Item(16, 4).Value = 1.339285714285714
$ws.Cells.Item(17, 3).Value = 378
$ws.Cells.Item(17, 4).Value = 2.116402116402116
$ws.Cells.Item(18, 3).Value = 328
$ws.Cells.Item(18, 4).Value = 2.439024390243902
$ws.Cells.Item(19, 3).Value = 342
$ws.Cells.Item(19, 4).Value = 2.923976608187134
$ws.Cells.Item(20, 2).Value = 7
$ws.Cells.Item(20, 3).Value = 299
$ws.Cells.Item(20, 4).Value = 2.341137123745819
$ws.Cells.Item(21, 3).Value = 234
$ws.Cells.Item(21, 4).Value = 1.70940170940171
$ws.Cells.Item(22, 2).Value = 7
$ws.Cells.Item(22, 3).Value = 205
$ws.Cells.Item(22, 4).Value = 3.414634146341464
$ws.Cells.Item(23, 3).Value = 160
$ws.Cells.Item(23, 4).Value = 2.5
$ws.Cells.Item(25, 3).Value = 271
$ws.Cells.Item(25, 4).Value = 0.7380073800738007
$ws.Cells.Item(27, 2).Value = 9
$ws.Cells.Item(27, 4).Value = 2.839116719242902
$ws.Cells.Item(28, 3).Value = 306
$ws.Cells.Item(28, 4).Value = 3.267973856209151
$ws.Cells.Item(29, 2).Value = 10
$ws.Cells.Item(29, 3).Value = 250
$ws.Cells.Item(29, 4).Value = 4
$ws.Cells.Item(30, 2).Value = 6
$ws.Cells.Item(30, 3).Value = 228
$ws.Cells.Item(30, 4).Value = 2.631578947368421
$ws.Cells.Item(31, 3).Value = 210
$ws.Cells.Item(31, 4).Value = 0.9523809523809524
$ws.Cells.Item(32, 3).Value = 120
$ws.Cells.Item(32, 4).Value = 0.8333333333333334
$ws.Cells.Item(33, 2).Value = 1
$ws.Cells.Item(33, 4).Value = 0.6369426751592357
$ws.Cells.Item(34, 3).Value = 243
$ws.Cells.Item(34, 4).Value = 2.469135802469136
$ws.Cells.Item(35, 2).Value = 6
$ws.Cells.Item(35, 3).Value = 230
$ws.Cells.Item(35, 4).Value = 2.608695652173913
$ws.Cells.Item(36, 3).Value = 243
$ws.Cells.Item(36, 4).Value = 4.938271604938271
$ws.Cells.Item(37, 2).Value = 11
$ws.Cells.Item(37, 3).Value = 245
$ws.Cells.Item(37, 4).Value = 4.489795918367347
$ws.Cells.Item(38, 2).Value = 8
$ws.Cells.Item(38, 3).Value = 275
$ws.Cells.Item(38, 4).Value = 2.909090909090909
$ws.Cells.Item(39, 2).Value = 12
$ws.Cells.Item(39, 3).Value = 224
$ws.Cells.Item(39, 4).Value = 5.357142857142857
$ws.Cells.Item(40, 3).Value = 220
$ws.Cells.Item(40, 4).Value = 3.181818181818182
$ws.Cells.Item(42, 3).Value = 289
$ws.Cells.Item(42, 4).Value = 3.114186851211072
$ws.Cells.Item(44, 2).Value = 6
$ws.Cells.Item(44, 3).Value = 250
$ws.Cells.Item(44, 4).Value = 2.4
$ws.Cells.Item(45, 2).Value = 8
$ws.Cells.Item(45, 4).Value = 3.389830508474576
$ws.Cells.Item(46, 3).Value = 266
$ws.Cells.Item(46, 4).Value = 5.263157894736842
$ws.Cells.Item(47, 3).Value = 242
$ws.Cells.Item(47, 4).Value = 4.545454545454546
$ws.Cells.Item(48, 3).Value = 280
$ws.Cells.Item(48, 4).Value = 4.285714285714286
$ws.Cells.Item(51, 2).Value = 22
$ws.Cells.Item(51, 3).Value = 300
$ws.Cells.Item(51, 4).Value = 7.333333333333333
$ws.Cells.Item(52, 3).Value = 292
$ws.Cells.Item(52, 4).Value = 3.767123287671233
$ws.Cells.Item(54, 2).Value = 8
$ws.Cells.Item(54, 3).Value = 235
$ws.Cells.Item(54, 4).Value = 3.404255319148936
$ws.Cells.Item(55, 2).Value = 5
$ws.Cells.Item(55, 3).Value = 249
$ws.Cells.Item(55, 4).Value = 2.008032128514056
$ws.Cells.Item(56, 3).Value = 295
$ws.Cells.Item(56, 4).Value = 5.084745762711865
$ws.Cells.Item(58, 3).Value = 270
$ws.Cells.Item(58, 4).Value = 2.962962962962963
$ws.Cells.Item(60, 2).Value = 13
$ws.Cells.Item(60, 4).Value = 4.436860068259386
$ws.Cells.Item(61, 2).Value = 9
$ws.Cells.Item(61, 3).Value = 260
$ws.Cells.Item(61, 4).Value = 3.461538461538462
$ws.Cells.Item(62, 2).Value = 10
$ws.Cells.Item(62, 3).Value = 247
$ws.Cells.Item(62, 4).Value = 4.048582995951417
$ws.Cells.Item(63, 3).Value = 238
$ws.Cells.Item(63, 4).Value = 3.361344537815126
$ws.Cells.Item(64, 3).Value = 283
$ws.Cells.Item(64, 4).Value = 2.120141342756184
$ws.Cells.Item(65, 3).Value = 207
$ws.Cells.Item(65, 4).Value = 2.415458937198068
$ws.Cells.Item(66, 2).Value = 8
$ws.Cells.Item(66, 3).Value = 223
$ws.Cells.Item(66, 4).Value = 3.587443946188341
$ws.Cells.Item(67, 3).Value = 258
$ws.Cells.Item(67, 4).Value = 3.10077519379845
$ws.Cells.Item(70, 3).Value = 235
$ws.Cells.Item(70, 4).Value = 7.234042553191489
$ws.Cells.Item(71, 3).Value = 240
$ws.Cells.Item(71, 4).Value = 8.333333333333332
$ws.Cells.Item(72, 3).Value = 190
$ws.Cells.Item(72, 4).Value = 8.421052631578947
$ws.Cells.Item(74, 2).Value = 23
$ws.Cells.Item(74, 3).Value = 266
$ws.Cells.Item(74, 4).Value = 8.646616541353383
$ws.Cells.Item(75, 2).Value = 15
$ws.Cells.Item(75, 3).Value = 259
$ws.Cells.Item(75, 4).Value = 5.791505791505791
$ws.Cells.Item(76, 2).Value = 11
$ws.Cells.Item(76, 3).Value = 229
$ws.Cells.Item(76, 4).Value = 4.803493449781659
$ws.Cells.Item(77, 2).Value = 17
$ws.Cells.Item(77, 3).Value = 256
$ws.Cells.Item(77, 4).Value = 6.640625
$ws.Cells.Item(79, 3).Value = 256
$ws.Cells.Item(79, 4).Value = 6.25
$ws.Cells.Item(80, 3).Value = 229
$ws.Cells.Item(80, 4).Value = 6.550218340611353
$ws.Cells.Item(81, 3).Value = 224
$ws.Cells.Item(81, 4).Value = 6.25
$ws.Cells.Item(82, 3).Value = 227
$ws.Cells.Item(82, 4).Value = 5.286343612334802
